# Group for LDAP integration
# Adds a "group type" dictionary entry (zero.group.type) under the
# existing zero.tabular table, plus its four concrete values
# (fixed / assignment / project / temp).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert a new row 6 (zero.tabular -> zero.group.type "组类型"),
#    pushing the existing zero.authority rows (old 6-9) down to 7-10.
#    Copy formatting from the row that is about to land below it so the
#    new row carries the same styles (s="8","11","5","9","9").
# ---------------------------------------------------------------
$ws.Rows.Item(6).Insert()
$ws.Range("A7:E7").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A6").Value = "b1ec6def-5adf-4f03-8b26-aadc4fee9e2e"
$ws.Range("B6").Value = "zero.tabular"
$ws.Range("C6").Value = 1010
$ws.Range("D6").Value = "组类型"
$ws.Range("E6").Value = "zero.group.type"

# ---------------------------------------------------------------
# 2. Append four new rows (11-14) describing the group-type values.
#    They reuse the same column layout: key, type, sort, name, code.
# ---------------------------------------------------------------
$ws.Range("A11").Value = "a594235a-b1df-4084-8ea3-58cf78f4361e"
$ws.Range("B11").Value = "zero.group.type"
$ws.Range("C11").Value = 1005
$ws.Range("D11").Value = "固定组"
$ws.Range("E11").Value = "fixed"

$ws.Range("A12").Value = "5906cf50-fe31-4319-8a94-9c8383157869"
$ws.Range("B12").Value = "zero.group.type"
$ws.Range("C12").Value = 1010
$ws.Range("D12").Value = "分派组"
$ws.Range("E12").Value = "assignment"

$ws.Range("A13").Value = "b88c89bc-7ea4-4047-a430-1e9da24c4ad8"
$ws.Range("B13").Value = "zero.group.type"
$ws.Range("C13").Value = 1015
$ws.Range("D13").Value = "项目组"
$ws.Range("E13").Value = "project"

$ws.Range("A14").Value = "a9ff7b0b-507e-47f1-a471-c8c62a04a8dc"
$ws.Range("B14").Value = "zero.group.type"
$ws.Range("C14").Value = 1020
$ws.Range("D14").Value = "临时组"
$ws.Range("E14").Value = "temp"

# Apply the same cell formatting (fill/border/font) used by the other
# data rows to the newly appended rows.
$ws.Range("A10:E10").Copy()
$ws.Range("A11:E14").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------
# 3. Match the saved selection left behind in the authored workbook.
# ---------------------------------------------------------------
$ws.Range("E15").Select()
